$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("NOMBRE PARTIDAS")

# Insert a new column before D for "Duracion (dias)" (shifts old D:J to E:K)
$ws1.Columns.Item(4).Insert()

# Row 1 header (new column D)
$ws1.Range("D1").Value = "Duración (días)"

# Row 2 (A:J) - full target state
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "Obra de prueba 1"
$ws1.Range("C2").Value = "Lima"
$ws1.Range("D2").Value = 60
$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = "2025-08-01"
$ws1.Range("E2").Style = "Normal"
$ws1.Range("F2").Value = "Limpieza del terreno"
$ws1.Range("G2").Value = "m²"
$ws1.Range("H2").Value = 100
$ws1.Range("I2").Value = 35
$ws1.Range("J2").Value = 3500

# Row 3 (A:J) - full target state
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Obra de prueba 1"
$ws1.Range("C3").Value = "Lima"
$ws1.Range("D3").Value = 60
$ws1.Range("E3").NumberFormat = "@"
$ws1.Range("E3").Value = "2025-08-01"
$ws1.Range("E3").Style = "Normal"
$ws1.Range("F3").Value = "Trazo y nivelación"
$ws1.Range("G3").Value = "m²"
$ws1.Range("H3").Value = 100
$ws1.Range("I3").Value = 80
$ws1.Range("J3").Value = 4000

# Row 4 (new)
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Obra de prueba 1"
$ws1.Range("C4").Value = "Lima"
$ws1.Range("D4").Value = 61
$ws1.Range("E4").NumberFormat = "@"
$ws1.Range("E4").Value = "2025-08-01"
$ws1.Range("E4").Style = "Normal"
$ws1.Range("F4").Value = "Excavación de zanjas"
$ws1.Range("G4").Value = "ml"
$ws1.Range("H4").Value = 300
$ws1.Range("I4").Value = 50
$ws1.Range("J4").Value = 15000

# Row 5 (new)
$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "Obra de prueba 1"
$ws1.Range("C5").Value = "Lima"
$ws1.Range("D5").Value = 62
$ws1.Range("E5").NumberFormat = "@"
$ws1.Range("E5").Value = "2025-08-01"
$ws1.Range("E5").Style = "Normal"
$ws1.Range("F5").Value = "Relleno compactado"
$ws1.Range("G5").Value = "m³"
$ws1.Range("H5").Value = 400
$ws1.Range("I5").Value = 45
$ws1.Range("J5").Value = 18000

# Clear the leftover column (old J shifted to K by the insert)
$ws1.Columns.Item(11).Clear()

$ws1.Columns.AutoFit()

$ws1.Range("J6").Select()
